# New topic: Drones in Agriculture
# Appends two new rows (16 = Predicted query row, 17 = Baseline row) to the
# results sheet, mirroring the layout of the existing Predicted/Baseline
# row pairs above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$queryText = @'

"precision agriculture" OR "geospatial analysis" OR "aerial imagery" OR "pesticide spraying" OR "crop health assessment" OR "drone technology" OR "autonomous drones" OR "crop insurance" OR "agroecology" OR "farm management software" OR "agricultural robotics" OR "field scouting" OR (drones in agriculture) OR "agricultural drone" OR 
(("yield estimation" OR "remote sensing" OR "variable rate application" OR "data analytics" OR "nutrient management" OR "land surveying" OR "unmanned aerial vehicle" OR "irrigation management" OR "climate monitoring" OR "sustainable farming" OR "farm productivity" OR "crop monitoring") AND (Drones OR Agriculture))

'@

# Use copy from an existing "Predicted"/"Baseline" label cell so the new
# label cells (column A) pick up the same bold/centered/bordered style
# (cellXfs index 1) instead of creating a brand-new duplicate style.
$ws.Range("A14").Copy($ws.Range("A16"))
$ws.Range("A15").Copy($ws.Range("A17"))

# Row 16 - Predicted
$ws.Range("A16").Value = "Predicted"
$ws.Range("B16").Value = $queryText
$ws.Range("C16").Value = 0.52
$ws.Range("D16").Value = 0.178
$ws.Range("E16").Value = 0.375
$ws.Range("F16").Value = 0.482
$ws.Range("G16").Value = 0.512

# The multi-line text just assigned makes the engine pin an explicit,
# auto-estimated row height (customHeight="1"); AutoFit re-measures the
# row and clears that explicit/custom flag again, matching the other
# "Predicted" rows that also hold multi-line query text.
$ws.Rows.Item(16).AutoFit()

# Row 17 - Baseline
$ws.Range("A17").Value = "Baseline"
$ws.Range("B17").Value = "Drones in Agriculture"
$ws.Range("C17").Value = 0.04
$ws.Range("D17").Value = 0.429
$ws.Range("E17").Value = 0.049
$ws.Range("F17").Value = 0.495
$ws.Range("G17").Value = 0.049
